$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")
$ws.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/documented-source"
